# Add a new column E that duplicates column C (rows 2-11), keeping the
# same per-row styling, then move the active selection to M4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy C2:C11 (values + style) into E2:E11, row by row, matching the
# original workbook's per-row style (rows 2-8 use style "3", rows 9-11
# use style "4" - Copy carries this over automatically).
$ws.Range("C2:C11").Copy($ws.Range("E2:E11")) | Out-Null

# Move the active cell / selection to M4, matching the saved workbook view.
$ws.Range("M4").Select() | Out-Null
